$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the coefficient/standard-error table with the new, more sensible values.
$ws.Range("B3").Value = "(0.076)"
$ws.Range("C3").Value = "(0.076)"

$ws.Range("B4").Value = "(0.089)"
$ws.Range("C4").Value = "(0.092)"

$ws.Range("B5").Value = "(0.109)"
$ws.Range("C5").Value = "(0.113)"

$ws.Range("B6").Value = "(0.128)"
$ws.Range("C6").Value = "(0.134)"

$ws.Range("B7").Value = "(0.149)"
$ws.Range("C7").Value = "(0.157)"

$ws.Range("B9").Value = "(0.107)"
$ws.Range("C9").Value = "(0.107)"

$ws.Range("B13").Value = "(0.032)"
$ws.Range("C13").Value = "(0.032)"

$ws.Range("B14").Value = "(0.135)"
$ws.Range("C14").Value = "(0.134)"

$ws.Range("C16").Value = "(0.056)"
$ws.Range("C17").Value = "(0.129)"
$ws.Range("C18").Value = "(0.057)"
$ws.Range("C19").Value = "(0.036)"
$ws.Range("C20").Value = "(0.024)"

$ws.Range("B21").Value = 8259
$ws.Range("C21").Value = 10222
